$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("104:104").Insert()
$ws.Range("A104").Value = "Noah"
$ws.Range("B104").Value = [char]0xea43
$ws.Range("C104").Value = 59971
$ws.Range("D104").Value = "ea43"

$win = $excel.ActiveWindow
$win.ScrollRow = 100
$win.ScrollColumn = 1
[void]$ws.Range("B105").Select()
